# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 18:52"

# --- Update country stats (row 4 = Estados Unidos) ---
$ws.Range("B4").Value = 930951
$ws.Range("C4").Value = 5719
$ws.Range("D4").Value = 110609
$ws.Range("E4").Value = 767428
$ws.Range("F4").Value = 15100
$ws.Range("G4").Value = 721
$ws.Range("H4").Value = 52914

# --- Row 10 = Iran ---
$ws.Range("B10").Value = 107773
$ws.Range("C10").Value = 2861
$ws.Range("D10").Value = 25582
$ws.Range("E10").Value = 79485
$ws.Range("F10").Value = 1782
$ws.Range("G10").Value = 106
$ws.Range("H10").Value = 2706

# --- Row 14 = Turquia ---
$ws.Range("B14").Value = 55224
$ws.Range("C14").Value = 2229
$ws.Range("E14").Value = 23807
$ws.Range("G14").Value = 92
$ws.Range("H14").Value = 3762

# --- Row 16 = Paises Bajos ---
$ws.Range("B16").Value = 44364
$ws.Range("C16").Value = 476
$ws.Range("E16").Value = 26545
$ws.Range("G16").Value = 48
$ws.Range("H16").Value = 2350

# --- Row 35 = Rumania ---
$ws.Range("E35").Value = 7144
$ws.Range("G35").Value = 34
$ws.Range("H35").Value = 601

# --- Row 47 = Malasia ---
$ws.Range("E47").Value = 4713
$ws.Range("G47").Value = 6
$ws.Range("H47").Value = 273

# --- Row 56 = Luxemburgo ---
$ws.Range("B56").Value = 3711
$ws.Range("C56").Value = 16
$ws.Range("E56").Value = 619

# --- Provincias / countries reorder near Nepal/Macao/Siria/Chad ---
# A new entry "Republica del Chad" is inserted right after Nepal (before Macao),
# shifting Macao's and Siria's existing figures down one row, with fresh figures
# for "Republica del Chad" on row 167.
$ws.Range("A167").Value = "Republica del Chad"
$ws.Range("B167").Value = 46
$ws.Range("C167").Value = 6
$ws.Range("D167").Value = 15
$ws.Range("E167").Value = 31
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 0

$ws.Range("A168").Value = "Macao"
$ws.Range("B168").Value = 45
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 27
$ws.Range("E168").Value = 18
$ws.Range("F168").Value = 1
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 0

$ws.Range("A169").Value = "Siria"
$ws.Range("B169").Value = 42
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 6
$ws.Range("E169").Value = 33
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 3
